$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for the columns that vary (D,I,J,K,L,M,O,P) across rows 2-24
# before overwriting, since the edit permutes rows.
$data = @{}
$data[2] = @{ "D" = 44980; "I" = "Primera"; "J" = 60; "K" = 7500; "L" = 8000; "M" = 7750; "O" = "Provincia de Diguillín"; "P" = 775 }
$data[3] = @{ "D" = 44782; "I" = "Primera"; "J" = 120; "K" = 8000; "L" = 9000; "M" = 8500; "O" = "Región Metropolitana"; "P" = 850 }
$data[4] = @{ "D" = 44211; "I" = "Primera"; "J" = 28; "K" = 8000; "L" = 8500; "M" = 8214; "O" = "Región Metropolitana"; "P" = 821 }
$data[5] = @{ "D" = 44838; "I" = "Primera"; "J" = 120; "K" = 6500; "L" = 7000; "M" = 6750; "O" = "Provincia de Diguillín"; "P" = 675 }
$data[6] = @{ "D" = 44817; "I" = "Primera"; "J" = 60; "K" = 7000; "L" = 7000; "M" = 7000; "O" = "Provincia de Diguillín"; "P" = 700 }
$data[7] = @{ "D" = 44817; "I" = "Segunda"; "J" = 60; "K" = 8000; "L" = 8000; "M" = 8000; "O" = "Provincia de Diguillín"; "P" = 800 }
$data[8] = @{ "D" = 44775; "I" = "Primera"; "J" = 60; "K" = 8000; "L" = 8000; "M" = 8000; "O" = "Región Metropolitana"; "P" = 800 }
$data[9] = @{ "D" = 44841; "I" = "Primera"; "J" = 60; "K" = 6500; "L" = 7000; "M" = 6750; "O" = "Provincia de Diguillín"; "P" = 675 }
$data[10] = @{ "D" = 44810; "I" = "Primera"; "J" = 60; "K" = 7000; "L" = 8000; "M" = 7500; "O" = "Provincia de Diguillín"; "P" = 750 }
$data[11] = @{ "D" = 44813; "I" = "Primera"; "J" = 120; "K" = 7000; "L" = 7500; "M" = 7250; "O" = "Provincia de Diguillín"; "P" = 725 }
$data[12] = @{ "D" = 44804; "I" = "Primera"; "J" = 80; "K" = 7000; "L" = 7500; "M" = 7250; "O" = "Provincia de Diguillín"; "P" = 725 }
$data[13] = @{ "D" = 44831; "I" = "Primera"; "J" = 60; "K" = 7000; "L" = 7500; "M" = 7250; "O" = "Provincia de Diguillín"; "P" = 725 }
$data[14] = @{ "D" = 44203; "I" = "Primera"; "J" = 27; "K" = 7000; "L" = 8000; "M" = 7556; "O" = "Región Metropolitana"; "P" = 756 }
$data[15] = @{ "D" = 44791; "I" = "Primera"; "J" = 100; "K" = 8500; "L" = 9000; "M" = 8750; "O" = "Región Metropolitana"; "P" = 875 }
$data[16] = @{ "D" = 44790; "I" = "Primera"; "J" = 60; "K" = 8500; "L" = 9000; "M" = 8750; "O" = "Región Metropolitana"; "P" = 875 }
$data[17] = @{ "D" = 44798; "I" = "Primera"; "J" = 80; "K" = 7000; "L" = 7000; "M" = 7000; "O" = "Provincia de Diguillín"; "P" = 700 }
$data[18] = @{ "D" = 44806; "I" = "Primera"; "J" = 120; "K" = 7000; "L" = 7500; "M" = 7250; "O" = "Provincia de Diguillín"; "P" = 725 }
$data[19] = @{ "D" = 44847; "I" = "Primera"; "J" = 100; "K" = 6500; "L" = 7000; "M" = 6750; "O" = "Provincia de Diguillín"; "P" = 675 }
$data[20] = @{ "D" = 44846; "I" = "Primera"; "J" = 100; "K" = 6500; "L" = 7000; "M" = 6750; "O" = "Provincia de Diguillín"; "P" = 675 }
$data[21] = @{ "D" = 44784; "I" = "Primera"; "J" = 100; "K" = 8000; "L" = 9000; "M" = 8500; "O" = "Región Metropolitana"; "P" = 850 }
$data[22] = @{ "D" = 44819; "I" = "Primera"; "J" = 100; "K" = 7000; "L" = 8000; "M" = 7500; "O" = "Provincia de Diguillín"; "P" = 750 }
$data[23] = @{ "D" = 44812; "I" = "Primera"; "J" = 60; "K" = 7000; "L" = 8000; "M" = 7500; "O" = "Provincia de Diguillín"; "P" = 750 }
$data[24] = @{ "D" = 44799; "I" = "Primera"; "J" = 60; "K" = 7000; "L" = 7000; "M" = 7000; "O" = "Provincia de Diguillín"; "P" = 700 }

# mapping: destination row -> source row (permutation derived from the edit)
$mapping = @{}
$mapping[2] = 4
$mapping[3] = 20
$mapping[4] = 10
$mapping[5] = 11
$mapping[6] = 24
$mapping[7] = 13
$mapping[8] = 14
$mapping[9] = 6
$mapping[10] = 7
$mapping[11] = 5
$mapping[12] = 2
$mapping[13] = 15
$mapping[14] = 8
$mapping[15] = 18
$mapping[16] = 21
$mapping[17] = 23
$mapping[18] = 19
$mapping[19] = 16
$mapping[20] = 22
$mapping[21] = 9
$mapping[22] = 17
$mapping[23] = 3
$mapping[24] = 12

# Apply the permutation: write each destination row using the snapshot of its source row
foreach ($destRow in 2..24) {
    $srcRow = $mapping[$destRow]
    $src = $data[$srcRow]
    $ws.Range("D" + $destRow).Value = $src["D"]
    $ws.Range("I" + $destRow).Value = $src["I"]
    $ws.Range("J" + $destRow).Value = $src["J"]
    $ws.Range("K" + $destRow).Value = $src["K"]
    $ws.Range("L" + $destRow).Value = $src["L"]
    $ws.Range("M" + $destRow).Value = $src["M"]
    $ws.Range("O" + $destRow).Value = $src["O"]
    $ws.Range("P" + $destRow).Value = $src["P"]
}